$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update cell values (columns C, D, E ; rows 2-9) ----
$ws.Range("C2").Value2 = -4.8629
$ws.Range("D2").Value2 = 0.6604
$ws.Range("E2").Value2 = 2.118

$ws.Range("C3").Value2 = -0.4464
$ws.Range("D3").Value2 = 0.4994
$ws.Range("E3").Value2 = 1.2062

$ws.Range("C4").Value2 = 0.4723
$ws.Range("D4").Value2 = 0.4027
$ws.Range("E4").Value2 = 0.9859

$ws.Range("C5").Value2 = 0.9219000000000001
$ws.Range("D5").Value2 = 0.1936
$ws.Range("E5").Value2 = 0.4756

$ws.Range("C6").Value2 = 0.6124000000000001
$ws.Range("D6").Value2 = 0.4823
$ws.Range("E6").Value2 = 1.2105

$ws.Range("C7").Value2 = 0.3245
$ws.Range("D7").Value2 = 0.6393
$ws.Range("E7").Value2 = 1.5981

$ws.Range("C8").Value2 = -0.0885
$ws.Range("D8").Value2 = 0.8250999999999999
$ws.Range("E8").Value2 = 2.0341

$ws.Range("C9").Value2 = -0.3501
$ws.Range("D9").Value2 = 0.9371
$ws.Range("E9").Value2 = 2.3226

# ---- Update the per-cell color-scale fills (RMSE column D, U column E) ----
# Color values below are plain RGB (alpha byte stripped) packed as 0xBBGGRR,
# the way the COM Interior.Color property expects them.

$ws.Range("D2").Interior.Color = 10279330   # 00A2D99C
$ws.Range("E2").Interior.Color = 14939879   # 00E7F6E3

$ws.Range("D3").Interior.Color = 6599248    # 0050B264
$ws.Range("E3").Interior.Color = 6401866    # 004AAF61

$ws.Range("D4").Interior.Color = 4887082    # 002A924A
$ws.Range("E4").Interior.Color = 4886825    # 0029914A

# D5 / E5 keep their previous fill color (0000441B) - no change needed.

$ws.Range("D6").Interior.Color = 6336070    # 0046AE60
$ws.Range("E6").Interior.Color = 6401866    # 004AAF61 (same shade as E3)

$ws.Range("D7").Interior.Color = 9754008    # 0098D594
$ws.Range("E7").Interior.Color = 9885339    # 009BD696

$ws.Range("D8").Interior.Color = 14349279   # 00DFF3DA
$ws.Range("E8").Interior.Color = 14283486   # 00DEF2D9

# D9 keeps its previous fill color (00F7FCF5) - no change needed.
$ws.Range("E9").Interior.Color = 16121079   # 00F7FCF5
